$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 7.824261088474345
$ws.Cells.Item(2, 4).Value = 6.807530062930876
$ws.Cells.Item(2, 5).Value = 9.773572391995119
$ws.Cells.Item(2, 6).Value = 37.03433097042019
$ws.Cells.Item(2, 7).Value = 3.679968072685935
$ws.Cells.Item(2, 11).Value = 15.17859874569648
$ws.Cells.Item(2, 13).Value = 16.42518571576159
$ws.Cells.Item(2, 14).Value = 21.02753937519212
$ws.Cells.Item(3, 2).Value = 7.754779508670498
$ws.Cells.Item(3, 4).Value = 6.826520897114252
$ws.Cells.Item(3, 5).Value = 9.644333653177222
$ws.Cells.Item(3, 6).Value = 36.68947545479829
$ws.Cells.Item(3, 7).Value = 3.683945753293896
$ws.Cells.Item(3, 11).Value = 14.81249415862314
$ws.Cells.Item(3, 13).Value = 16.19663982001731
$ws.Cells.Item(3, 14).Value = 21.06748923603365
$ws.Cells.Item(4, 2).Value = 7.713725099500479
$ws.Cells.Item(4, 4).Value = 6.838861953500285
$ws.Cells.Item(4, 5).Value = 9.567241304595193
$ws.Cells.Item(4, 6).Value = 36.48671385748267
$ws.Cells.Item(4, 7).Value = 3.686512532372289
$ws.Cells.Item(4, 11).Value = 14.58829678414626
$ws.Cells.Item(4, 13).Value = 16.05994050078433
$ws.Cells.Item(4, 14).Value = 21.09389555744258
$ws.Cells.Item(5, 2).Value = 7.697417151736683
$ws.Cells.Item(5, 4).Value = 6.844062076652903
$ws.Cells.Item(5, 5).Value = 9.536429453851163
$ws.Cells.Item(5, 6).Value = 36.40641926017597
$ws.Cells.Item(5, 7).Value = 3.687589940230948
$ws.Cells.Item(5, 11).Value = 14.49723385547207
$ws.Cells.Item(5, 13).Value = 16.0052154255271
$ws.Cells.Item(5, 14).Value = 21.10512775087132
$ws.Cells.Item(6, 2).Value = 7.694735225711356
$ws.Cells.Item(6, 4).Value = 6.844935880972384
$ws.Cells.Item(6, 5).Value = 9.531350650025647
$ws.Cells.Item(6, 6).Value = 36.39322933546831
$ws.Cells.Item(6, 7).Value = 3.687770744654925
$ws.Cells.Item(6, 11).Value = 14.4821354708636
$ws.Cells.Item(6, 13).Value = 15.99618957783029
$ws.Cells.Item(6, 14).Value = 21.10702130318248
$ws.Cells.Item(7, 2).Value = 7.713503432963323
$ws.Cells.Item(7, 4).Value = 6.838931391799365
$ws.Cells.Item(7, 5).Value = 9.566823274272286
$ws.Cells.Item(7, 6).Value = 36.48562143899758
$ws.Cells.Item(7, 7).Value = 3.686526935269079
$ws.Cells.Item(7, 11).Value = 14.58706726430222
$ws.Cells.Item(7, 13).Value = 16.05919839836572
$ws.Cells.Item(7, 14).Value = 21.0940451308779
$ws.Cells.Item(8, 2).Value = 7.799981121019071
$ws.Cells.Item(8, 4).Value = 6.813936778539734
$ws.Cells.Item(8, 5).Value = 9.728563332758172
$ws.Cells.Item(8, 6).Value = 36.91360672949188
$ws.Cells.Item(8, 7).Value = 3.681313820016489
$ws.Cells.Item(8, 11).Value = 15.05233189209553
$ws.Cells.Item(8, 13).Value = 16.34567167246924
$ws.Cells.Item(8, 14).Value = 21.04092418186383
$ws.Cells.Item(9, 2).Value = 7.981499944431238
$ws.Cells.Item(9, 4).Value = 6.770327448857846
$ws.Cells.Item(9, 5).Value = 10.06205463582595
$ws.Cells.Item(9, 6).Value = 37.82065797419023
$ws.Cells.Item(9, 7).Value = 3.672072741958657
$ws.Cells.Item(9, 11).Value = 15.96304309170919
$ws.Cells.Item(9, 13).Value = 16.93313810882891
$ws.Cells.Item(9, 14).Value = 20.95167276518701
$ws.Cells.Item(10, 2).Value = 8.120999676585964
$ws.Cells.Item(10, 4).Value = 6.741589498099418
$ws.Cells.Item(10, 5).Value = 10.31479164372174
$ws.Cells.Item(10, 6).Value = 38.52342841112491
$ws.Cells.Item(10, 7).Value = 3.665873791429857
$ws.Cells.Item(10, 11).Value = 16.62302219406168
$ws.Cells.Item(10, 13).Value = 17.37618885259293
$ws.Cells.Item(10, 14).Value = 20.89523063602732
$ws.Cells.Item(11, 2).Value = 8.185543599374371
$ws.Cells.Item(11, 4).Value = 6.729234638338383
$ws.Cells.Item(11, 5).Value = 10.43096358550105
$ws.Cells.Item(11, 6).Value = 38.84991110645272
$ws.Cells.Item(11, 7).Value = 3.663180213622772
$ws.Cells.Item(11, 11).Value = 16.91971600099581
$ws.Cells.Item(11, 13).Value = 17.5793278770345
$ws.Cells.Item(11, 14).Value = 20.87154450931041
$ws.Cells.Item(12, 2).Value = 8.210119342180551
$ws.Cells.Item(12, 4).Value = 6.724659661249071
$ws.Cells.Item(12, 5).Value = 10.47508742093925
$ws.Cells.Item(12, 6).Value = 38.97442057039875
$ws.Cells.Item(12, 7).Value = 3.662178262163886
$ws.Cells.Item(12, 11).Value = 17.03143879160166
$ws.Cells.Item(12, 13).Value = 17.65640623860917
$ws.Cells.Item(12, 14).Value = 20.86286205069808
$ws.Cells.Item(13, 2).Value = 8.204820885242366
$ws.Cells.Item(13, 4).Value = 6.725640354957439
$ws.Cells.Item(13, 5).Value = 10.46557933116755
$ws.Cells.Item(13, 6).Value = 38.94756778310885
$ws.Cells.Item(13, 7).Value = 3.66239324958464
$ws.Cells.Item(13, 11).Value = 17.00740707694299
$ws.Cells.Item(13, 13).Value = 17.63980036225148
$ws.Cells.Item(13, 14).Value = 20.86471919695094
$ws.Cells.Item(14, 2).Value = 8.187562891755363
$ws.Cells.Item(14, 4).Value = 6.728856175518418
$ws.Cells.Item(14, 5).Value = 10.43459126419106
$ws.Cells.Item(14, 6).Value = 38.86013755677807
$ws.Cells.Item(14, 7).Value = 3.663097421462362
$ws.Cells.Item(14, 11).Value = 16.92892076153928
$ws.Cells.Item(14, 13).Value = 17.58566649769426
$ws.Cells.Item(14, 14).Value = 20.87082444210315
$ws.Cells.Item(15, 2).Value = 8.177008727491287
$ws.Cells.Item(15, 4).Value = 6.730839450731627
$ws.Cells.Item(15, 5).Value = 10.41562617880231
$ws.Cells.Item(15, 6).Value = 38.80669535339025
$ws.Cells.Item(15, 7).Value = 3.663531093840433
$ws.Cells.Item(15, 11).Value = 16.88076032357721
$ws.Cells.Item(15, 13).Value = 17.55252579983181
$ws.Cells.Item(15, 14).Value = 20.87460147772488
$ws.Cells.Item(16, 2).Value = 8.116801560287914
$ws.Cells.Item(16, 4).Value = 6.742411391165852
$ws.Cells.Item(16, 5).Value = 10.30722021055161
$ws.Cells.Item(16, 6).Value = 38.50222020222625
$ws.Cells.Item(16, 7).Value = 3.666052355256504
$ws.Cells.Item(16, 11).Value = 16.60355093843689
$ws.Cells.Item(16, 13).Value = 17.36293884561492
$ws.Cells.Item(16, 14).Value = 20.89681868974791
$ws.Cells.Item(17, 2).Value = 8.080128885203433
$ws.Cells.Item(17, 4).Value = 6.749694527781469
$ws.Cells.Item(17, 5).Value = 10.24099511999719
$ws.Cells.Item(17, 6).Value = 38.31710357200885
$ws.Cells.Item(17, 7).Value = 3.667631343469044
$ws.Cells.Item(17, 11).Value = 16.43249827033672
$ws.Cells.Item(17, 13).Value = 17.24698754770997
$ws.Cells.Item(17, 14).Value = 20.91095843275219
$ws.Cells.Item(18, 2).Value = 8.059139140677592
$ws.Cells.Item(18, 4).Value = 6.753951176423341
$ws.Cells.Item(18, 5).Value = 10.20301961334821
$ws.Cells.Item(18, 6).Value = 38.21127439653858
$ws.Cells.Item(18, 7).Value = 3.668551436207888
$ws.Cells.Item(18, 11).Value = 16.33378731967882
$ws.Cells.Item(18, 13).Value = 17.18044978500919
$ws.Cells.Item(18, 14).Value = 20.91927846510772
$ws.Cells.Item(19, 2).Value = 8.052050792436823
$ws.Cells.Item(19, 4).Value = 6.755404003817205
$ws.Cells.Item(19, 5).Value = 10.19018282723535
$ws.Cells.Item(19, 6).Value = 38.17555622324547
$ws.Cells.Item(19, 7).Value = 3.668865011433363
$ws.Cells.Item(19, 11).Value = 16.30031331724851
$ws.Cells.Item(19, 13).Value = 17.15795006344242
$ws.Cells.Item(19, 14).Value = 20.92212761576971
$ws.Cells.Item(20, 2).Value = 8.084022208267614
$ws.Cells.Item(20, 4).Value = 6.748912229204751
$ws.Cells.Item(20, 5).Value = 10.2480332448833
$ws.Cells.Item(20, 6).Value = 38.33674350776909
$ws.Cells.Item(20, 7).Value = 3.667462026707439
$ws.Cells.Item(20, 11).Value = 16.4507417577154
$ws.Cells.Item(20, 13).Value = 17.25931533441353
$ws.Cells.Item(20, 14).Value = 20.90943384969283
$ws.Cells.Item(21, 2).Value = 8.192628508825196
$ws.Cells.Item(21, 4).Value = 6.727908798781843
$ws.Cells.Item(21, 5).Value = 10.44368994776622
$ws.Cells.Item(21, 6).Value = 38.88579490907908
$ws.Cells.Item(21, 7).Value = 3.662890100263865
$ws.Cells.Item(21, 11).Value = 16.95199207475486
$ws.Cells.Item(21, 13).Value = 17.60156333614824
$ws.Cells.Item(21, 14).Value = 20.86902338862001
$ws.Cells.Item(22, 2).Value = 8.264382130278618
$ws.Cells.Item(22, 4).Value = 6.714785493130083
$ws.Cells.Item(22, 5).Value = 10.57231418983361
$ws.Cells.Item(22, 6).Value = 39.2497003691071
$ws.Cells.Item(22, 7).Value = 3.660007223235817
$ws.Cells.Item(22, 11).Value = 17.27587385972845
$ws.Cells.Item(22, 13).Value = 17.82610634650003
$ws.Cells.Item(22, 14).Value = 20.84428586824229
$ws.Cells.Item(23, 2).Value = 8.226022278856032
$ws.Cells.Item(23, 4).Value = 6.721734320542998
$ws.Cells.Item(23, 5).Value = 10.50360952124231
$ws.Cells.Item(23, 6).Value = 39.05504625532055
$ws.Cells.Item(23, 7).Value = 3.661536288712426
$ws.Cells.Item(23, 11).Value = 17.10338944024417
$ws.Cells.Item(23, 13).Value = 17.70620881080123
$ws.Cells.Item(23, 14).Value = 20.85733537783314
$ws.Cells.Item(24, 2).Value = 8.082261742902478
$ws.Cells.Item(24, 4).Value = 6.749265690007879
$ws.Cells.Item(24, 5).Value = 10.24485099972701
$ws.Cells.Item(24, 6).Value = 38.32786242631518
$ws.Cells.Item(24, 7).Value = 3.667538536458865
$ws.Cells.Item(24, 11).Value = 16.44249502183322
$ws.Cells.Item(24, 13).Value = 17.25374154953214
$ws.Cells.Item(24, 14).Value = 20.91012251912928
$ws.Cells.Item(25, 2).Value = 7.931232124759576
$ws.Cells.Item(25, 4).Value = 6.781545654092461
$ws.Cells.Item(25, 5).Value = 9.970317190210915
$ws.Cells.Item(25, 6).Value = 37.56851873706272
$ws.Cells.Item(25, 7).Value = 3.674468423480928
$ws.Cells.Item(25, 11).Value = 15.71771427411785
$ws.Cells.Item(25, 13).Value = 16.77189670410579
$ws.Cells.Item(25, 14).Value = 20.97421684317742
